# Update field "形態" (data type) for CreateDate and LastUpdate rows from
# DATE to TIMESTAMP on the DBD sheet, then leave the selection on the last
# cell edited (D14), matching the interactive edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

$ws.Range("D12").Value = "TIMESTAMP"
$ws.Range("D14").Value = "TIMESTAMP"

$ws.Range("D14").Select()
